$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated designators on BOM
$ws.Range("I10").Value = "X$1"
$ws.Range("I11").Value = "Mates to X$1"
$ws.Range("I6").Value = "IC1Mx"

# Added a price for the TO220 mounting kit
$ws.Range("H9").Value = 1.57

# Restore the last-used selection
$null = $ws.Range("I7").Select()
